$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new profit entry for 2025-11-22 as row 97, right after the
# existing last row (96, which holds 11/21/2025).
# The date is entered with a leading apostrophe so Excel stores it as
# literal text ("11/22/2025"), matching how the other date cells in this
# column are stored, instead of being auto-converted to a date serial
# number. ClearFormats() then drops the quote-prefix formatting flag that
# gets attached to the cell so it keeps the sheet's default (unstyled)
# look, same as its neighboring data cells.
$ws.Range("A97").Value = "'11/22/2025"
$ws.Range("A97").ClearFormats()

$ws.Range("B97").Value = 7565.55
